# Edit script: remove the "Price" column, shift FirstName/LastName left,
# and append two new client rows (Jet/Hassan/Baraka and Mars/Gabriel/Alfaro).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column B ("Price") entirely - this shifts C (FirstName) and D (LastName)
# left into B and C, matching the new 3-column layout (ItemName, FirstName, LastName).
$ws.Range("B1").EntireColumn.Delete()

# Add the two new client rows under the existing data.
$ws.Range("A7").Value = "Jet"
$ws.Range("B7").Value = "Hassan"
$ws.Range("C7").Value = "Baraka"

$ws.Range("A8").Value = "Mars"
$ws.Range("B8").Value = "Gabriel"
$ws.Range("C8").Value = "Alfaro"

# Update the active selection as reflected in the saved workbook.
$ws.Range("E10").Select()
